$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Quarterly rollover: shift columns E:N left by one quarter and append the new quarter ---

# Header label rows (quarter captions), columns E(5)..N(14)
$quarterLabels = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)
$headerRows = @(8, 16, 26, 35, 43, 52)
foreach ($r in $headerRows) {
    for ($i = 0; $i -lt 10; $i++) {
        $ws.Cells.Item($r, 5 + $i).Value = $quarterLabels[$i]
    }
}

# Data rows: new values for columns E(5)..N(14) after the rollover
$row10 = @("-", "-", "-", 0, 0, 0, 0, 0, 0, 0)
$row11 = @(1950409450, 1518329270, 1482260280, 1747522533, 1602405913, 1461425094, 1647975740, 1633774747, 1592843945, 1383769998)
$row12 = @(1950409450, 1518329270, 1482260280, 1747522533, 1602405913, 1461425094, 1647975740, 1633774747, 1592843945, 1383769998)
$row18 = @(0, -1594000, "-", "-", "-", "-", "-", "-", "-", "-")
$row19 = @(0, 0, "-", 0, 0, 0, 0, 0, 0, 0)
$row20 = @(1625067530, 1380570350, 1554315700, "-", 1632643800, 1364693200, 1630766000, 1627185900, 1591083500, 1384813600)
$row21 = @("-", "-", "-", 1715422250, "-", "-", "-", "-", "-", "-")
$row22 = @(1625067530, 1378976350, 1554315700, 1715422250, 1632643800, 1364693200, 1630766000, 1627185900, 1591083500, 1384813600)
$row28 = @(0, -259, "-", "-", "-", "-", "-", "-", "-", "-")
$row29 = @(0, 0, "-", 0, 0, 0, 0, 0, 0, 0)
$row30 = @(608479, 521243, 748084, 839829, 1020360, 848430, 1442520, 1440442, 1404910, 1228423)
$row31 = @(608479, 520984, 748084, 839829, 1020360, 848430, 1442520, 1440442, 1404910, 1228423)
$row37 = @("-", 162, "-", "-", "-", "-", "-", "-", "-", "-")
$row38 = @("-", "-", "-", "-", "-", "-", "-", "-", "-", "-")
$row39 = @(375, 378, 481, 490, 625, 622, 885, 885, 883, 887)
$row45 = @(110, 186, "-", "-", "-", "-", "-", "-", "-", "-")
$row46 = @(0, 0, "-", 0, 0, 0, 0, 0, 0, 0)
$row47 = @(-259171, -121304, -357694, -442826, -438367, -381522, -549891, -476995, -614024, -429191)
$row48 = @(-259061, -121118, -357694, -442826, -438367, -381522, -549891, -476995, -614024, -429191)
$row54 = @(110, -73, "-", "-", "-", "-", "-", "-", "-", "-")
$row55 = @(0, 0, "-", 0, 0, 0, 0, 0, 0, 0)
$row56 = @(349175, 399939, 390390, 397001, 581993, 466908, 892629, 963447, 790886, 799232)
$row57 = @(349285, 399866, 390390, 397001, 581993, 466908, 892629, 963447, 790886, 799232)
$dataRows = @{
    10 = $row10
    11 = $row11
    12 = $row12
    18 = $row18
    19 = $row19
    20 = $row20
    21 = $row21
    22 = $row22
    28 = $row28
    29 = $row29
    30 = $row30
    31 = $row31
    37 = $row37
    38 = $row38
    39 = $row39
    45 = $row45
    46 = $row46
    47 = $row47
    48 = $row48
    54 = $row54
    55 = $row55
    56 = $row56
    57 = $row57
}
foreach ($r in $dataRows.Keys) {
    $vals = $dataRows[$r]
    for ($i = 0; $i -lt 10; $i++) {
        $ws.Cells.Item([int]$r, 5 + $i).Value = $vals[$i]
    }
}
